$d = $word.ActiveDocument

# --- Change 1: BIOS settings sentence ---
$d.Content.Find.Execute(
    "(Execute Disable on Intel systems; No Execute on AMD) in your BIOS settings.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Execute Disable on Intel systems; No Execute on AMD) are enabled in your BIOS settings.",
    2)

# --- Change 2 & 3: Question 1 merge ---
$d.Content.Find.Execute(
    "Question 1: With regards to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Question 1: With regards to ",
    2)

$d.Content.Find.Execute(
    ", what does indeterminate mean? Why is this indicator indeterminate?  (12 points)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", what does indeterminate mean? Why is this indicator indeterminate?  (12 points)",
    2)
